$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.71471856507537
$ws.Range("D2").Value = 9.445537643816754
$ws.Range("E2").Value = 15.86926128448355
$ws.Range("F2").Value = 37.62063430798291
$ws.Range("G2").Value = 3.684714637473324
$ws.Range("J2").Value = 11.77034564934962
$ws.Range("K2").Value = 9.484535063601168
$ws.Range("L2").Value = 8.863847523803175
$ws.Range("O2").Value = 28.32498482443728
$ws.Range("B3").Value = 16.61125011914996
$ws.Range("D3").Value = 9.451859020083331
$ws.Range("E3").Value = 15.91153393480013
$ws.Range("F3").Value = 37.71730575754095
$ws.Range("G3").Value = 3.686687749209808
$ws.Range("J3").Value = 11.79721946560308
$ws.Range("K3").Value = 9.229893326226893
$ws.Range("L3").Value = 8.81983037974919
$ws.Range("O3").Value = 28.41091671872402
$ws.Range("B4").Value = 16.55031555506875
$ws.Range("D4").Value = 9.456863924032866
$ws.Range("E4").Value = 15.93922085408454
$ws.Range("F4").Value = 37.78440960724957
$ws.Range("G4").Value = 3.687964181293232
$ws.Range("J4").Value = 11.81459967364422
$ws.Range("K4").Value = 9.071013392392608
$ws.Range("L4").Value = 8.793713857152131
$ws.Range("O4").Value = 28.46914371865661
$ws.Range("B5").Value = 16.52615780228728
$ws.Range("D5").Value = 9.459186625627352
$ws.Range("E5").Value = 15.95093962524227
$ws.Range("F5").Value = 37.81369993344803
$ws.Range("G5").Value = 3.688500716914712
$ws.Range("J5").Value = 11.82190405258698
$ws.Range("K5").Value = 9.005726698227214
$ws.Range("L5").Value = 8.783307293418313
$ws.Range("O5").Value = 28.49424384532512
$ws.Range("B6").Value = 16.5221877000414
$ws.Range("D6").Value = 9.459589428461683
$ws.Range("E6").Value = 15.95291188186337
$ws.Range("F6").Value = 37.81868094408778
$ws.Range("G6").Value = 3.688590799013074
$ws.Range("J6").Value = 11.82313035394534
$ws.Range("K6").Value = 8.994856023019741
$ws.Range("L6").Value = 8.78159375979731
$ws.Range("O6").Value = 28.49849451861067
$ws.Range("B7").Value = 16.54998700098157
$ws.Range("D7").Value = 9.45689410146907
$ws.Range("E7").Value = 15.93937713064186
$ws.Range("F7").Value = 37.78479675610444
$ws.Range("G7").Value = 3.687971350813648
$ws.Range("J7").Value = 11.81469728414715
$ws.Range("K7").Value = 9.070134978742562
$ws.Range("L7").Value = 8.793572545125318
$ws.Range("O7").Value = 28.46947667459214
$ws.Range("B8").Value = 16.67851679904454
$ws.Range("D8").Value = 9.447484454820497
$ws.Range("E8").Value = 15.88347801619955
$ws.Range("F8").Value = 37.65235688856502
$ws.Range("G8").Value = 3.685381520275082
$ws.Range("J8").Value = 11.77942956308864
$ws.Range("K8").Value = 9.397314181674588
$ws.Range("L8").Value = 8.848485939183034
$ws.Range("O8").Value = 28.35347885453244
$ws.Range("B9").Value = 16.95019468745924
$ws.Range("D9").Value = 9.437918019530542
$ws.Range("E9").Value = 15.7875628305825
$ws.Range("F9").Value = 37.45424406827942
$ws.Range("G9").Value = 3.680815756367495
$ws.Range("J9").Value = 11.71722002461352
$ws.Range("K9").Value = 10.01500212996287
$ws.Range("L9").Value = 8.963062911188825
$ws.Range("O9").Value = 28.16944929221953
$ws.Range("B10").Value = 17.16041899668197
$ws.Range("D10").Value = 9.436265950642543
$ws.Range("E10").Value = 15.72539966578729
$ws.Range("F10").Value = 37.34640953846034
$ws.Range("G10").Value = 3.677770694999479
$ws.Range("J10").Value = 11.67571150247202
$ws.Range("K10").Value = 10.4495197421763
$ws.Range("L10").Value = 9.050987646895441
$ws.Range("O10").Value = 28.06083628840675
$ws.Range("B11").Value = 17.25806132844628
$ws.Range("D11").Value = 9.436672167523204
$ws.Range("E11").Value = 15.69891345370418
$ws.Range("F11").Value = 37.30557245419568
$ws.Range("G11").Value = 3.676451908272782
$ws.Range("J11").Value = 11.65773093458005
$ws.Range("K11").Value = 10.64212838357241
$ws.Range("L11").Value = 9.091698591099636
$ws.Range("O11").Value = 28.01722164535341
$ws.Range("B12").Value = 17.29529889351317
$ws.Range("D12").Value = 9.436991554243843
$ws.Range("E12").Value = 15.6891407631967
$ws.Range("F12").Value = 37.29129200298416
$ws.Range("G12").Value = 3.675962017596288
$ws.Range("J12").Value = 11.65105120578836
$ws.Range("K12").Value = 10.71427274813083
$ws.Range("L12").Value = 9.107208522586713
$ws.Range("O12").Value = 28.00154070700481
$ws.Range("B13").Value = 17.28726786764908
$ws.Range("D13").Value = 9.436915420497931
$ws.Range("E13").Value = 15.69123406414722
$ws.Range("F13").Value = 37.29431488149689
$ws.Range("G13").Value = 3.676067102313302
$ws.Range("J13").Value = 11.65248407111009
$ws.Range("K13").Value = 10.69877138492024
$ws.Range("L13").Value = 9.103864164390776
$ws.Range("O13").Value = 28.004880714087
$ws.Range("B14").Value = 17.26111977704278
$ws.Range("D14").Value = 9.436695131046763
$ws.Range("E14").Value = 15.69810430029612
$ws.Range("F14").Value = 37.30437386016962
$ws.Range("G14").Value = 3.676411414439681
$ws.Range("J14").Value = 11.657178804694
$ws.Range("K14").Value = 10.64807997662919
$ws.Range("L14").Value = 9.092972779161865
$ws.Range("O14").Value = 28.01591482091198
$ws.Range("B15").Value = 17.24513671236767
$ws.Range("D15").Value = 9.436581730038634
$ws.Range("E15").Value = 15.70234597659856
$ws.Range("F15").Value = 37.31068947707892
$ws.Range("G15").Value = 3.676623552036518
$ws.Range("J15").Value = 11.66007126338246
$ws.Range("K15").Value = 10.61692500160244
$ws.Range("L15").Value = 9.086313408070565
$ws.Range("O15").Value = 28.02278232080105
$ws.Range("B16").Value = 17.15407603202288
$ws.Range("D16").Value = 9.436262622787162
$ws.Range("E16").Value = 15.72716660909674
$ws.Range("F16").Value = 37.34924387325121
$ws.Range("G16").Value = 3.677858213477937
$ws.Range("J16").Value = 11.67690467429577
$ws.Range("K16").Value = 10.43682485740312
$ws.Range("L16").Value = 9.048340680736411
$ws.Range("O16").Value = 28.06380335526101
$ws.Range("B17").Value = 17.098709368382
$ws.Range("D17").Value = 9.436362750367621
$ws.Range("E17").Value = 15.74285180986994
$ws.Range("F17").Value = 37.37500196327464
$ws.Range("G17").Value = 3.678632618341636
$ws.Range("J17").Value = 11.68746201254459
$ws.Range("K17").Value = 10.32499595753458
$ws.Range("L17").Value = 9.025222019428277
$ws.Range("O17").Value = 28.09045373455767
$ws.Range("B18").Value = 17.06705475706789
$ws.Range("D18").Value = 9.436529372843401
$ws.Range("E18").Value = 15.75204224356607
$ws.Range("F18").Value = 37.39059069603634
$ws.Range("G18").Value = 3.679084290735402
$ws.Range("J18").Value = 11.69361923573836
$ws.Range("K18").Value = 10.26020259053678
$ws.Range("L18").Value = 9.011992488357734
$ws.Range("O18").Value = 28.10632758155556
$ws.Range("B19").Value = 17.05637062692294
$ws.Range("D19").Value = 9.436604541829743
$ws.Range("E19").Value = 15.75518296418591
$ws.Range("F19").Value = 37.39600153130119
$ws.Range("G19").Value = 3.679238295018394
$ws.Range("J19").Value = 11.69571856837127
$ws.Range("K19").Value = 10.23818563567515
$ws.Range("L19").Value = 9.007525098733051
$ws.Range("O19").Value = 28.11179578389115
$ws.Range("B20").Value = 17.10458368517386
$ws.Range("D20").Value = 9.43634081375064
$ws.Range("E20").Value = 15.74116463615197
$ws.Range("F20").Value = 37.37217991493306
$ws.Range("G20").Value = 3.678549534573375
$ws.Range("J20").Value = 11.68632938092365
$ws.Range("K20").Value = 10.33694971305826
$ws.Range("L20").Value = 9.027676097469014
$ws.Range("O20").Value = 28.08756031089765
$ws.Range("B21").Value = 17.26879319326472
$ws.Range("D21").Value = 9.436755349740038
$ws.Range("E21").Value = 15.69607937386177
$ws.Range("F21").Value = 37.30138715564365
$ws.Range("G21").Value = 3.676310024023999
$ws.Range("J21").Value = 11.65579634778338
$ws.Range("K21").Value = 10.66299126709501
$ws.Range("L21").Value = 9.096169375837363
$ws.Range("O21").Value = 28.01265116128661
$ws.Range("B22").Value = 17.37763130878027
$ws.Range("D22").Value = 9.437990823227544
$ws.Range("E22").Value = 15.66811161688694
$ws.Range("F22").Value = 37.26201998801207
$ws.Range("G22").Value = 3.67490175672341
$ws.Range("J22").Value = 11.63659360895259
$ws.Range("K22").Value = 10.87143374958704
$ws.Range("L22").Value = 9.141474880051893
$ws.Range("O22").Value = 27.9685613589433
$ws.Range("B23").Value = 17.31941234065226
$ws.Range("D23").Value = 9.437243490376245
$ws.Range("E23").Value = 15.68290167200375
$ws.Range("F23").Value = 37.28239906352433
$ws.Range("G23").Value = 3.675648323391648
$ws.Range("J23").Value = 11.64677381621332
$ws.Range("K23").Value = 10.76062909302471
$ws.Range("L23").Value = 9.117247974763739
$ws.Range("O23").Value = 27.9916469580503
$ws.Range("B24").Value = 17.10192735500197
$ws.Range("D24").Value = 9.436350391537276
$ws.Range("E24").Value = 15.74192686908522
$ws.Range("F24").Value = 37.37345333327796
$ws.Range("G24").Value = 3.678587076629184
$ws.Range("J24").Value = 11.68684117055552
$ws.Range("K24").Value = 10.33154697704201
$ws.Range("L24").Value = 9.02656641593191
$ws.Range("O24").Value = 28.0888667078138
$ws.Range("B25").Value = 16.87473286044637
$ws.Range("D25").Value = 9.43955869167843
$ws.Range("E25").Value = 15.81204842680163
$ws.Range("F25").Value = 37.50122572101429
$ws.Range("G25").Value = 3.681996347287273
$ws.Range("J25").Value = 11.73330945418084
$ws.Range("K25").Value = 9.850972480878639
$ws.Range("L25").Value = 8.931375203369427
$ws.Range("O25").Value = 28.21457163275858
